$d = $word.ActiveDocument
$d.Content.Find.Execute("Online Server", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Database", 2)
